$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODO_TEAM")

# Row 11: clear the now-obsolete E11 cell and update D11's status text
# ("Ouvert" -> "Corrigé -> to be test"), reusing/recycling the shared
# string that used to hold "TO BE TEST ".
$ws.Range("E11").ClearContents()
$ws.Range("D11").Value = "Corrigé -> to be test"

# Row 12: status moved from "Ouvert" to "Corrigé"
$ws.Range("D12").Value = "Corrigé"

# Row 8 (A8): apply a time number format (h:mm) to the task cell
$ws.Range("A8").NumberFormat = "h:mm"

# Update the active selection to A8
$ws.Range("A8").Select()
